$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Meetings")
Write-Host $ws.Name
